$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the values that need to shift down from rows 14-17 to rows 20-23
$vSecReg      = $ws.Range("A14").Value()   # mtoSecReg[n]
$vSecRegMnd   = $ws.Range("A15").Value()   # mtoSecReg[mnd][n]
$vSecRegPct   = $ws.Range("A16").Value()   # mtoSecRegPct[n]
$vSecRegGraf  = $ws.Range("A17").Value()   # mtoSecRegGraf

# Clear the old rows 14-17 before repopulating (rows 16/19 end up blank)
$ws.Range("A14:C17").ClearContents()

# New dictionary entries inserted at rows 14, 15, 17, 18
$ws.Range("A14").Value = "mtoInsReg[n]"
$ws.Range("A15").Value = "mtoInsReg[mnd][n]"
$ws.Range("A17").Value = "mtoRubReg[n]"
$ws.Range("A18").Value = "mtoRubReg[mnd][n]"

# Annotation added next to the new dictionary entry
$ws.Range("B14").Value = "QUERY???"

# The four "mtoSecReg*" rows move down to rows 20-23
$ws.Range("A20").Value = $vSecReg
$ws.Range("A21").Value = $vSecRegMnd
$ws.Range("A22").Value = $vSecRegPct
$ws.Range("A23").Value = $vSecRegGraf

# Annotation added next to the relocated mtoSecReg[n] entry
$ws.Range("B20").Value = "M???"

# Update the view: scroll so row 7 is the first visible row, select B20
$ws.Range("B20").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
